# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary block ---
# VALOR MORA total
$ws.Cells.Item(11, 5).Value2 = 54858

# Cant. Trabajadores / Cant. Periodos
$ws.Cells.Item(13, 3).Value2 = 2
$ws.Cells.Item(13, 6).Value2 = 2

# --- Row 16: first worker (now LILIEN TAINA BARRIOS OLIVO) ---
$ws.Cells.Item(16, 2).Value2 = "CC"
$ws.Cells.Item(16, 3).Value2 = "1049929498"
$ws.Cells.Item(16, 4).Value2 = "LILIEN TAINA BARRIOS OLIVO"
$ws.Cells.Item(16, 5).Value2 = "2105"
$ws.Cells.Item(16, 6).Value2 = 9306
$ws.Cells.Item(16, 7).Value2 = 2235014

# --- Row 17 becomes the new last data row (MILFRE LILIANA MOSQUERA IBARGUEN) ---
# Pull the "closing border" formatting that currently lives on row 20 down onto row 17
# before the old rows get removed, so the final table keeps its bottom border.
$ws.Range("B20:J20").Copy() | Out-Null
$ws.Range("B17:J17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(17, 2).Value2 = "CC"
$ws.Cells.Item(17, 3).Value2 = "35851218"
$ws.Cells.Item(17, 4).Value2 = "MILFRE LILIANA MOSQUERA IBARGUEN"
$ws.Cells.Item(17, 5).Value2 = "2506"
$ws.Cells.Item(17, 6).Value2 = 45552
$ws.Cells.Item(17, 7).Value2 = 1423500
$ws.Cells.Item(17, 8).Value2 = ""
$ws.Cells.Item(17, 9).Value2 = ""
$ws.Cells.Item(17, 10).Value2 = ""

# --- Remove the now-obsolete rows 18, 19 and 20 ---
$ws.Rows("18:20").Delete()

# --- Resize columns to fit the new, shorter content (values taken from the
# workbook's bestFit recalculation after the data refresh) ---
$ws.Columns("B").ColumnWidth = 16.90625
$ws.Columns("C").ColumnWidth = 10.81640625
$ws.Columns("D").ColumnWidth = 32.54296875
$ws.Columns("E").ColumnWidth = 12.7265625
$ws.Columns("F").ColumnWidth = 9.453125
$ws.Columns("G").ColumnWidth = 13.453125
$ws.Columns("H").ColumnWidth = 17.90625
$ws.Columns("I").ColumnWidth = 16.81640625
$ws.Columns("J").ColumnWidth = 14.1796875

# Nudge the logo image back in line with the new, narrower column widths
$shp = $ws.Shapes.Item(1)
$shp.Left = $shp.Left - 16.319448818897638
$shp2 = $ws.Shapes.Item(1)
$shp2.Width = 76.81889763779527
$shp2.Height = 48.188976377952756
